$d = $word.ActiveDocument

# --- Fix 1: byline run-split correction (typo fix in names) ---
# Originally split as: "Taylor " | [spellStart]"Okonek"[spellEnd] | " & Charlie Wolock"
# Target split is:      "Taylor Okonek & Charlie " | [spellStart]"Wolock"[spellEnd]
# Work on narrow, single-run matches so Word doesn't re-flow the whole
# paragraph into one run and so proofErr tags stay put.

# 1) Remove the trailing " & Charlie Wolock" run entirely.
$d.Content.Find.Execute(" & Charlie Wolock", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2)

# 2) Turn the spell-checked "Okonek" run into "Wolock" (do this before
#    growing the first run so "Okonek" still only matches that one run).
$d.Content.Find.Execute("Okonek", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Wolock", 2)

# 3) Expand the leading "Taylor " run to include the rest of the byline.
$d.Content.Find.Execute("Taylor ", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Taylor Okonek & Charlie ", 2)

# --- Fix 2: mislabeled list item ---
# The document has two "Early MRI" / "Surgery" variable lists -- one for
# the "early MRI and surgery" question, one for the "early MRI and
# opioid prescription" question. Only the list item belonging to the
# latter question is mislabeled "Surgery" and needs to read
# "Opioid prescription " instead. Find it by walking paragraphs and
# remembering the most recent question stem we saw.
$sawOpioidStem = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    $t = $para.Range.Text
    if ($t -match "early MRI and opioid prescription") {
        $sawOpioidStem = $true
    } elseif ($t -match "early MRI and surgery") {
        $sawOpioidStem = $false
    }
    if ($sawOpioidStem -and $t -eq "Surgery`r") {
        $para.Range.Find.Execute("Surgery", $true, $false, $false, $false, $false, `
                                  $true, 1, $false, "Opioid prescription ", 2)
        $sawOpioidStem = $false
    }
}
